# Add 7 new rows (590-596) of landscaping data to Sheet1, extend the
# shared ABS() formula in column F down through the new rows, and move
# the active selection to reflect the new bottom of the data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string columns are written as their literal text; Excel's
# COM layer re-interns them into xl/sharedStrings.xml on save.
$newRows = @(
    @{R=590; A=45871; B="Flowering";     C="Large";  D=56; E=78; G=0; H=0;    I="No"; J=2; K="Bright";  L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7},
    @{R=591; A=45871; B="Nonflowering";  C="Medium"; D=56; E=78; G=0; H=0;    I="No"; J=3; K="Bright";  L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7},
    @{R=592; A=45871; B="Nonflowering";  C="Small";  D=56; E=78; G=0; H=0.1;  I="No"; J=3; K="Neutral"; L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7},
    @{R=593; A=45871; B="Nonflowering";  C="Medium"; D=56; E=78; G=0; H=0.15; I="No"; J=3; K="Neutral"; L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7},
    @{R=594; A=45871; B="Nonflowering";  C="Medium"; D=56; E=78; G=0; H=0.1;  I="No"; J=3; K="Bright";  L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7},
    @{R=595; A=45871; B="Nonflowering";  C="Large";  D=56; E=78; G=0; H=0.25; I="No"; J=4; K="Neutral"; L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7},
    @{R=596; A=45871; B="Tree";          C="Medium"; D=56; E=78; G=0; H=0.5;  I="No"; J=1; K="Neutral"; L=8; M=0.42; N=51; O=30.29; P=9; Q=0.35; R2=9.9; S=46; T=7}
)

foreach ($row in $newRows) {
    $r = $row.R

    # Column A carries the same short-date number format as the rows
    # above it; copy that formatting over before writing the value so
    # the existing date style (s="1") is reused instead of a new one
    # being minted.
    $ws.Range("A$r").Value2 = $row.A
    $ws.Range("A589").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, "B").Value2 = $row.B
    $ws.Cells.Item($r, "C").Value2 = $row.C
    $ws.Cells.Item($r, "D").Value2 = $row.D
    $ws.Cells.Item($r, "E").Value2 = $row.E
    # F (Temp_Diff) is filled in afterwards as a single shared formula.
    $ws.Cells.Item($r, "G").Value2 = $row.G
    $ws.Cells.Item($r, "H").Value2 = $row.H
    $ws.Cells.Item($r, "I").Value2 = $row.I
    $ws.Cells.Item($r, "J").Value2 = $row.J
    $ws.Cells.Item($r, "K").Value2 = $row.K
    $ws.Cells.Item($r, "L").Value2 = $row.L
    $ws.Cells.Item($r, "M").Value2 = $row.M
    $ws.Cells.Item($r, "N").Value2 = $row.N
    $ws.Cells.Item($r, "O").Value2 = $row.O
    $ws.Cells.Item($r, "P").Value2 = $row.P
    $ws.Cells.Item($r, "Q").Value2 = $row.Q
    $ws.Cells.Item($r, "R").Value2 = $row.R2
    $ws.Cells.Item($r, "S").Value2 = $row.S
    $ws.Cells.Item($r, "T").Value2 = $row.T
}

# Extend the ABS(D-E) formula down through the new rows in one shot so
# it's written back out as a single shared formula (same mechanism the
# original file used for F543:F589).
$ws.Range("F590:F596").Formula = "=ABS(D590-E590)"

# Match the author's final viewport/selection: scrolled near the new
# bottom of the sheet, with the last filled column (T) selected for the
# newly-added rows.
$excel.Goto($ws.Range("A574"), $true)
$ws.Range("T590:T596").Select()
